$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Mos Eisley Banthas (7)'
$ws.Range("B2").Value = 'Team Blah (4)'

$ws.Range("A7").Value = 'Please Hammer don''t Huerter (3)'
$ws.Range("B7").Value = 'Team Winchester (6)'
$ws.Range("D7").Value = 'Team Winchester'

$ws.Range("A15").Value = 'Team James (3)'
$ws.Range("B15").Value = 'The A**L Embiids (7)'

$ws.Range("A16").Value = 'Edward St Easybeats (4)'
$ws.Range("B16").Value = 'VICTOR-IOUS (6)'
$ws.Range("D16").Value = 'VICTOR-IOUS'

$ws.Range("A22").Value = 'Team Cooke (6)'
$ws.Range("B22").Value = 'The David Cahill''s (4)'
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 'Team Cooke'

$ws.Range("A24").Value = 'Balls Deep (7)'
$ws.Range("B24").Value = 'Team luka slo (3)'

$ws.Range("A25").Value = 'Browntown (3)'
$ws.Range("C25").Value = 2

$ws.Range("A31").Value = 'Boswell Pioneers (7)'
$ws.Range("B31").Value = 'misq (3)'

$ws.Range("A41").Value = 'Scottie Pippings (5)'
$ws.Range("B41").Value = 'Zions Illegitimate Child (4)'
$ws.Range("C41").Value = 2

$ws.Range("A44").Value = 'Meep Meep (6)'
$ws.Range("B44").Value = 'Ol'' Dirty Baskets (5)'
$ws.Range("C44").Value = 0

$ws.Range("A45").Value = 'Bulls (6)'
$ws.Range("B45").Value = 'Team Hewish (5)'

$ws.Range("A46").Value = 'Bank Town Squids (6)'
$ws.Range("B46").Value = 'You Don''t Mess with the Ant-Man (4)'
$ws.Range("D46").Value = 'Bank Town Squids'

$ws.Range("A51").Value = 'Brisbane Noble’s Nobs (4)'
$ws.Range("C51").Value = 1

$ws.Range("A52").Value = 'Ja Raffe (7)'
$ws.Range("B52").Value = 'Team Mercer (3)'

$ws.Range("A54").Value = 'Big Bam Theory (5)'
$ws.Range("B54").Value = 'Triple Towers (6)'
$ws.Range("C54").Value = 0

$ws.Range("B62").Value = 'Moon Shooters (5)'
$ws.Range("C62").Value = 2

$ws.Range("A64").Value = 'Boomtown Bulls (4)'
$ws.Range("B64").Value = 'Shanghai Sharks (6)'
$ws.Range("D64").Value = 'Shanghai Sharks'

$ws.Range("A65").Value = 'Boston Ballerz (3)'
$ws.Range("B65").Value = 'Duncans Donuts (7)'
